$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update y_corrSteps (E), y_nrSteps (G), alienID (H) for the affected trial rows.
# Sheet row -> (E, G, H) new values
$updates = @{
    4  = @{ E = 6; G = 3; H = 13 }
    8  = @{ E = 6; G = 3; H = 13 }
    16 = @{ E = 7; G = 3; H = 13 }
    18 = @{ E = 6; G = 3; H = 13 }
    23 = @{ E = 5; G = 3; H = 13 }
    27 = @{ E = 7; G = 3; H = 13 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
    $ws.Range("H$row").Value = $vals.H
}
